$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Start Time" values (column F) for every scheduled class move forward by
# exactly one day (the underlying date serial increases by 1) while keeping
# the same time-of-day. Set the exact target serial values directly so the
# stored doubles match full precision.
$newValues = @{
    2  = 3.4166666666666665
    3  = 3.4166666666666665
    4  = 3.5
    5  = 3.5
    6  = 3.5833333333333335
    8  = 3.4166666666666665
    9  = 3.4166666666666665
    10 = 3.5
    11 = 3.5833333333333335
    12 = 3.5833333333333335
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value2 = $newValues[$row]
}
